# Updated cryptos list with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row,
# and for a handful of rows the Coin/Link (B/C) pair was also re-ranked.
# Price cells are stored as text (matching the source data's inline-string
# cells), so NumberFormat is forced to "@" before writing any value that
# would otherwise be auto-parsed as a number by Excel (which would drop
# meaningful trailing zeros, e.g. "6.60" -> 6.6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.450.97"
$ws.Range("E2").Value = "  +3.55%  "

$ws.Range("D3").Value = "2.258.53"
$ws.Range("E3").Value = "  +1.99%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.10"
$ws.Range("E5").Value = "  +2.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.78"
$ws.Range("E6").Value = "  +4.74%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.521"
$ws.Range("E7").Value = "  +2.23%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.486"
$ws.Range("E9").Value = "  +2.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "31.98"
$ws.Range("E10").Value = "  +6.52%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.99"
$ws.Range("E11").Value = "  +7.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0795"
$ws.Range("E12").Value = "  +2.31%  "

$ws.Range("E13").Value = "  +2.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.60"
$ws.Range("E14").Value = "  +2.22%  "

$ws.Range("D15").Value = "2.602.46"
$ws.Range("E15").Value = "  +1.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.15"
$ws.Range("E16").Value = "  +2.67%  "

$ws.Range("D17").Value = "2.241.20"
$ws.Range("E17").Value = "  +0.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.752"
$ws.Range("E18").Value = "  +2.90%  "

$ws.Range("D19").Value = "41.350.76"
$ws.Range("E19").Value = "  +3.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.98"
$ws.Range("E20").Value = "  +6.60%  "

$ws.Range("D21").Value = "0.0₃0904"
$ws.Range("E21").Value = "  +2.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.85"
$ws.Range("E22").Value = "  +1.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.86"
$ws.Range("E23").Value = "  +2.16%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.32"
$ws.Range("E24").Value = "  +2.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.56"
$ws.Range("E25").Value = "  +3.99%  "

$ws.Range("E26").Value = "  -0.89%  "

$ws.Range("E27").Value = "  +3.75%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.83"
$ws.Range("E28").Value = "  +5.39%  "

$ws.Range("E29").Value = "  -1.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.56"
$ws.Range("E30").Value = "  +3.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.04"
$ws.Range("E31").Value = "  +2.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.46"
$ws.Range("E32").Value = "  +4.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.17"
$ws.Range("E34").Value = "  +5.44%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0734"
$ws.Range("E35").Value = "  +2.67%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.03"
$ws.Range("E36").Value = "  +5.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.37"
$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.105"
$ws.Range("E38").Value = "  +8.22%  "

$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.51"
$ws.Range("E39").Value = "  +6.15%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.115"
$ws.Range("E40").Value = "  +2.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.80"
$ws.Range("E41").Value = "  +7.50%  "

$ws.Range("E42").Value = "  +3.39%  "

$ws.Range("D43").Value = "2.074.94"
$ws.Range("E43").Value = "  -2.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.20"
$ws.Range("E44").Value = "  +12.66%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0278"
$ws.Range("E45").Value = "  +4.37%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.23"
$ws.Range("E46").Value = "  +4.84%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.96"
$ws.Range("E47").Value = "  +10.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.04"
$ws.Range("E48").Value = "  -3.19%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.476.28"
$ws.Range("E49").Value = "  +1.69%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.52"
$ws.Range("E50").Value = "  +4.18%  "

$ws.Range("E51").Value = "  +5.13%  "

